# Swap the data values between row 2 and row 3 for columns A, B, E, F, G, I, M
# (the other columns already hold identical values in both rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "M")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}

# Column I ("Antal") is stored as text ("1"/"2"), not a number, so force
# text formatting for the assignment and then drop back to the default
# (unstyled) cell style so no stray number-format style is introduced.
$i2 = $ws.Range("I2")
$i3 = $ws.Range("I3")

$iv2 = $i2.Value2
$iv3 = $i3.Value2

$i2.NumberFormat = "@"
$i2.Value2 = $iv3
$i2.Style = "Normal"

$i3.NumberFormat = "@"
$i3.Value2 = $iv2
$i3.Style = "Normal"
